$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the header text of A1: "Projected_Revenute" -> "Projected_Revenue"
$ws.Range("A1").Value = "Projected_Revenue"

# Clear the stale B14 selection/active-cell state left over from editing,
# moving the cursor back to A1 (the neutral/default position).
$ws.Range("A1").Select()
